$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 754 (pushes existing rows 754:795 down to 755:796)
$ws.Rows.Item(754).EntireRow.Insert()

# Column A holds dates formatted as plain text in this sheet (not real
# Excel dates), so force text formatting before writing the value, then
# strip the format again so the cell ends up with no explicit style -
# matching the rest of the data rows.
$ws.Cells.Item(754, 1).NumberFormat = "@"
$ws.Cells.Item(754, 1).Value = "2026/01/31"
$ws.Cells.Item(754, 1).ClearFormats()

$ws.Cells.Item(754, 2).Value = "土"
$ws.Cells.Item(754, 3).Value = 8
$ws.Cells.Item(754, 4).Value = 201
